# Apply updated cryptocurrency price/volume data to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "29.301.03"
$ws.Range("E2").Value = "  +3.03%  "
Set-TextValue "D3" "1.903.51"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  -1.30%  "
Set-TextValue "D5" "315.50"
$ws.Range("E5").Value = "  -0.56%  "
Set-TextValue "D6" "1.003"
$ws.Range("E6").Value = "  -1.47%  "
Set-TextValue "D7" "0.5138"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  -0.49%  "
Set-TextValue "D9" "0.08468"
$ws.Range("E9").Value = "  +0.30%  "
Set-TextValue "D10" "42.50"
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("E11").Value = "  +0.74%  "
Set-TextValue "D12" "6.265"
$ws.Range("E12").Value = "  +0.13%  "
Set-TextValue "D13" "1.903.41"
$ws.Range("E13").Value = "  +1.79%  "
Set-TextValue "D14" "20.60"
$ws.Range("E14").Value = "  +0.55%  "
Set-TextValue "D15" "7.361"
$ws.Range("E15").Value = "  +1.59%  "
Set-TextValue "D16" "1.005"
$ws.Range("E16").Value = "  -1.32%  "
Set-TextValue "D17" "93.27"
$ws.Range("E17").Value = "  +2.31%  "
Set-TextValue "D18" "0.00001108"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("E20").Value = "  +1.12%  "
Set-TextValue "D21" "1.003"
$ws.Range("E21").Value = "  -1.33%  "
Set-TextValue "D22" "6.042"
$ws.Range("E22").Value = "  +1.49%  "
Set-TextValue "D23" "29.297.35"
$ws.Range("E23").Value = "  +2.93%  "
Set-TextValue "D25" "2.219"
$ws.Range("E25").Value = "  -2.43%  "
Set-TextValue "D26" "2.119.24"
$ws.Range("E26").Value = "  +1.72%  "
Set-TextValue "D27" "160.21"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  +0.62%  "
Set-TextValue "D29" "2.458"
$ws.Range("E29").Value = "  +3.98%  "
Set-TextValue "D30" "127.94"
$ws.Range("E30").Value = "  +0.94%  "
Set-TextValue "D31" "0.1050"
$ws.Range("E31").Value = "  -0.54%  "
Set-TextValue "D32" "1.060"
$ws.Range("E32").Value = "  +1.54%  "
Set-TextValue "D33" "6.091"
$ws.Range("E33").Value = "  +5.53%  "
Set-TextValue "D34" "3.653"
$ws.Range("E34").Value = "  +0.77%  "
Set-TextValue "D35" "0.02480"
$ws.Range("E35").Value = "  +1.95%  "
Set-TextValue "D36" "0.06599"
$ws.Range("E36").Value = "  +1.73%  "
Set-TextValue "D37" "9.136"
$ws.Range("E37").Value = "  +2.52%  "
Set-TextValue "D38" "0.2200"
$ws.Range("E38").Value = "  +1.14%  "
Set-TextValue "D39" "1.240"
$ws.Range("E39").Value = "  +4.39%  "
Set-TextValue "D40" "5.134"
$ws.Range("E40").Value = "  +2.54%  "
Set-TextValue "D41" "0.6528"
$ws.Range("E41").Value = "  +1.89%  "
Set-TextValue "D42" "1.236"
$ws.Range("E42").Value = "  -2.12%  "
Set-TextValue "D43" "11.30"
$ws.Range("E43").Value = "  +0.65%  "
Set-TextValue "D44" "0.6058"
$ws.Range("E44").Value = "  +0.12%  "
Set-TextValue "D45" "13.17"
$ws.Range("E45").Value = "  +1.51%  "
Set-TextValue "D46" "3.680"
$ws.Range("E46").Value = "  -1.13%  "
Set-TextValue "D47" "2.059"
$ws.Range("E47").Value = "  +3.43%  "
$ws.Range("E48").Value = "  +2.11%  "
Set-TextValue "D49" "123.05"
$ws.Range("E49").Value = "  +0.68%  "
Set-TextValue "D50" "1.176"
$ws.Range("E50").Value = "  -2.16%  "
Set-TextValue "D51" "77.88"
$ws.Range("E51").Value = "  +1.53%  "
